# Update metrics values in columns B:Q for rows 2 through 26 with the
# new values produced by the retraining ("atualizado todo o treinamento
# para o novo lm"). All data rows share the same updated metric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B..Q (same set applied to every data row).
# Values in scientific notation are parsed via [double] cast from string
# because the PowerShell-style parser here doesn't accept bare "1e-05"
# literals.
$newValues = @(
    [double]"0.9999677858346675",
    [double]"0.9990493920370204",
    [double]"0.9999584699111314",
    [double]"0.9999381311800511",
    [double]"0.9999573860719446",
    [double]"3.007050929946918e-05",
    [double]"0.0008873508065744985",
    [double]"4.729565515708677e-05",
    [double]"7.275692831993349e-05",
    [double]"6.002629165744963e-05",
    [double]"0.0003164627591308915",
    [double]"0.005483658386466938",
    [double]"1.000026659998896",
    [double]"0.005717109076716589",
    [double]"126.8239312479565",
    [double]"191.4243499659711"
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - 2]
    }
}
